# "random data & output update"
# Adds two new commit rows (31 and 32) to the Arena commits table, extends
# the Total(h) formula to include them, and updates the window's scroll /
# selection state to match the author's final view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: new commit entry
$ws.Range("C31").Value = "#include & constructors update & source added"
$ws.Range("G31").Value = 3.5

# Row 32: new commit entry
$ws.Range("C32").Value = "random data & output update"
$ws.Range("G32").Value = 2.5

# Extend the Total(h) sum so it picks up the two new rows
$ws.Range("G39").Formula = "=SUM(G4:G32)"

# Match the saved view: scrolled down with G33 selected
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G33").Select()
